# Apply the "fixing frontend design in attendance table" edit.
# The sheet shrinks from a 6-column x 7-row attendance table to a
# 5-column x 3-row table: column E (the "2025-03-13" date column) is
# dropped, the old column F slides left into column E, the D1 header
# text is replaced, rows 4-7 (extra students) are removed, and the
# remaining attendance text values are replaced with numeric codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole column E (shifts old F into E, removing the old
# "2025-03-13" column entirely), then delete rows 4 through 7 (the
# extra student rows that are no longer present).
$ws.Range("E1:E7").EntireColumn.Delete()
$ws.Range("A4:F7").EntireRow.Delete()

# Update header text for the remaining columns.
$ws.Range("D1").Value = "2025-03-13 - asdasdadaasd"
$ws.Range("E1").Value = "23-01-2025 - HW1"

# Update attendance values for the two remaining student rows with the
# new numeric coding.
$ws.Range("D2").Value = "0"
$ws.Range("E2").Value = "1"

$ws.Range("D3").Value = "0"
$ws.Range("E3").Value = "0"

# Restore the custom column widths (they shift left by one relative to
# the removed column, and columns B-E get new widths).
$ws.Columns.Item(1).Width = 14.83203125
$ws.Columns.Item(2).Width = 14.83203125
$ws.Columns.Item(3).Width = 19.83203125
$ws.Columns.Item(4).Width = 27.83203125
$ws.Columns.Item(5).Width = 18.83203125
